# Auto-evalutation de l'avancement du projet
#
# Highlights each checklist item of the evaluation sheet with a colored
# status (red = not started, green = done, gold/orange = in progress,
# amber/theme = partially done) and adds a few review comments in
# column C, plus clears the (no longer meaningful) percentage formula
# for the first section and tweaks the selection / page setup.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Red = not started yet
# ---------------------------------------------------------------------
$redCells = @("B3", "B4", "B5", "B6", "B7", "B8", "B9", "B21", "B29", "B30")
foreach ($addr in $redCells) {
    $ws.Range($addr).Interior.Color = 255
}

# ---------------------------------------------------------------------
# Green = done
# ---------------------------------------------------------------------
$greenCells = @("B12", "B13", "B14", "B15", "B25", "B31", "B36")
foreach ($addr in $greenCells) {
    $ws.Range($addr).Interior.Color = 5287936
}

# ---------------------------------------------------------------------
# Orange/gold = in progress
# ---------------------------------------------------------------------
$orangeCells = @("B32", "B33", "B34", "B35", "B37", `
                  "B41", "B42", "B43", "B44", "B45", "B46", "B47", "B48")
foreach ($addr in $orangeCells) {
    $ws.Range($addr).Interior.Color = 49407
}

# ---------------------------------------------------------------------
# Theme accent (amber) = partially done
# ---------------------------------------------------------------------
$themeCells = @("B16", "B17", "B18", "B19", "B20", "B22", "B23", "B24")
foreach ($addr in $themeCells) {
    $ws.Range($addr).Interior.ThemeColor = 10
}

# ---------------------------------------------------------------------
# The completion percentage is no longer tracked for the first section
# ---------------------------------------------------------------------
$ws.Range("C9").Value = $null

# ---------------------------------------------------------------------
# Review comments added in column C
# ---------------------------------------------------------------------
$ws.Range("C21").Value = "!!!!!!!!!!!!!!!!!!!!!!!!!!!!!!"
$ws.Range("C21").Font.Bold = $true
$ws.Range("C21").Font.Color = 255

$ws.Range("C23").Value = "(OK pour la création, faire pour la modification)"
$ws.Range("C24").Value = "??? Voir avec le prof… la politique"
$ws.Range("C34").Value = "Vérifier que c'est correct"
$ws.Range("C37").Value = "A vérifier"
$ws.Range("C41").Value = "Il n'y a pas encore les EndUsers"

# ---------------------------------------------------------------------
# Misc view / print tweaks
# ---------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("B31").Select()
